$d = $word.ActiveDocument

# 1) Merge the three "Life " / "Cycle" / " Method:" runs into a single run.
$d.Content.Find.Execute("Life Cycle Method:", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Life Cycle Method:", 2)

# 2) Merge the two "...they" / " use props." runs into a single run.
$d.Content.Find.Execute("they use props.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "they use props.", 2)

# 3) Append three new bulleted ("ListParagraph", numId 1) paragraphs at the
#    very end of the document body, right after "Here Children is the children."
$end = $d.Content.End
$r = $d.Range($end, $end)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
         '<w:pPr>' + `
           '<w:pStyle w:val="ListParagraph"/>' + `
           '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
         '</w:pPr>' + `
         '<w:r><w:t>React uses Syntheic Events.</w:t></w:r>' + `
       '</w:p>' + `
       '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
         '<w:pPr>' + `
           '<w:pStyle w:val="ListParagraph"/>' + `
           '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
         '</w:pPr>' + `
         '<w:r><w:t>Destructruing in React:</w:t></w:r>' + `
         '<w:r><w:tab/><w:t>const  {monsters} = this.state.monsters</w:t></w:r>' + `
       '</w:p>' + `
       '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
         '<w:pPr>' + `
           '<w:pStyle w:val="ListParagraph"/>' + `
           '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
         '</w:pPr>' + `
         '<w:r><w:t>Functional Components unlike class component component donot have access to state(no access to constructor), and life cycle methods.</w:t></w:r>' + `
       '</w:p>'

$r.InsertXML($xml)
